# Updates the "cryptos" price/volume table (Sheet1) for the latest GitHub
# Actions scrape: refreshed Price (col D) and Volume(1h) (col E) figures,
# plus a few coins that swapped ranking positions (rows 23/24, 40/41/42,
# 49/50), which carries their Coin name + Link along with them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D stores prices as plain text (e.g. "212.17"), even when the text
# looks like a number. Pre-format the cells whose new price text would
# otherwise be auto-parsed as a numeric value by Excel, so they keep being
# stored as text -- matching every other cell in this column.
foreach ($addr in @('D5','D6','D14','D16','D19','D20','D23','D24','D25','D30','D31','D39','D40','D41','D42','D44','D45','D47','D49','D50')) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '26.273.16'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.589.39'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('D5').Value = '213.32'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').Value = '0.501'
$ws.Range('E6').Value = '  +0.67%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '1.812.76'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = '1.591.15'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').Value = '4.05'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').Value = '64.43'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').Value = '26.281.69'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('D19').Value = '7.47'
$ws.Range('E19').Value = '  +2.03%  '
$ws.Range('D20').Value = '214.09'
$ws.Range('E20').Value = '  +3.04%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = '2.15'
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('B24').Value = 'Avalanche'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D24').Value = '8.95'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('D25').Value = '145.22'
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').Value = '0.0500'
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('D34').Value = '1.343.77'
$ws.Range('E34').Value = '  +5.17%  '
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('E37').Value = '  -2.45%  '
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('D39').Value = '0.816'
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '5.78'
$ws.Range('E40').Value = '  +4.11%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '1.01'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('D44').Value = '0.763'
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('D45').Value = '61.85'
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('D46').Value = '1.724.88'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').Value = '87.41'
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('E48').Value = '  -4.35%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0504'
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0979'
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('E51').Value = '  -0.31%  '
